# CDS_SPARSE_2014_2015.xlsx rewrite: cohort graduation-rate columns restructured.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): columns C, E, G..N get new labels (shared strings
#    are managed automatically by the engine - unused ones get pruned, new
#    ones appended).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value  = "stafford-loan"
$ws.Cells.Item(1, 5).Value  = "any-aid"
$ws.Cells.Item(1, 7).Value  = "exemptions"
$ws.Cells.Item(1, 8).Value  = "final"
$ws.Cells.Item(1, 9).Value  = "between"
$ws.Cells.Item(1, 10).Value = "within"
$ws.Cells.Item(1, 11).Value = "4 year"
$ws.Cells.Item(1, 12).Value = "5 year"
$ws.Cells.Item(1, 13).Value = "6 year"
$ws.Cells.Item(1, 14).Value = "initial"

# ---------------------------------------------------------------------------
# 2. Row 5: N5 flips from 0 to 1 (initial -> final cohort reclassified).
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 14).Value = 1

# ---------------------------------------------------------------------------
# 3. Rows 14-25: columns I..M (between, within, 4 year, 5 year, 6 year) are
#    recoded from a single "stage" indicator to the new split scheme.
# ---------------------------------------------------------------------------
$rowsBetweenWithin4 = @(14, 15, 16, 17)
foreach ($r in $rowsBetweenWithin4) {
    $ws.Cells.Item($r, 9).Value  = 0
    $ws.Cells.Item($r, 10).Value = 1
    $ws.Cells.Item($r, 11).Value = 1
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
}

$rows5 = @(18, 19, 20, 21)
foreach ($r in $rows5) {
    $ws.Cells.Item($r, 9).Value  = 1
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 1
    $ws.Cells.Item($r, 12).Value = 1
    $ws.Cells.Item($r, 13).Value = 0
}

$rows6 = @(22, 23, 24, 25)
foreach ($r in $rows6) {
    $ws.Cells.Item($r, 9).Value  = 1
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 1
    $ws.Cells.Item($r, 13).Value = 1
}

# ---------------------------------------------------------------------------
# 4. Drop the trailing summary rows (26-34); the sheet now ends at row 25.
# ---------------------------------------------------------------------------
$ws.Rows("26:34").Delete()

# ---------------------------------------------------------------------------
# 5. Strip the red/plain-font override (style index 2) from columns E and H
#    - cells fall back to the default (automatic/theme) font.
# ---------------------------------------------------------------------------
$ws.Range("E1:E25").ClearFormats()
$ws.Range("H1:H25").ClearFormats()

# ---------------------------------------------------------------------------
# 6. Column widths: only B:C (15 chars) and G (11.33 chars) keep an explicit
#    width now.
# ---------------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 14.17
$ws.Columns("G:G").ColumnWidth = 10.5

# ---------------------------------------------------------------------------
# 7. Selection cosmetic update.
# ---------------------------------------------------------------------------
$ws.Range("E5").Select()
